$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text without Excel coercing
# numeric-looking strings (e.g. "241.54") into a Number, and without
# leaving the cell style changed (captures/restores the original Style).
function Set-TextValue($range, [string]$text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") '41.942.93'
$ws.Range("E2").Value = '  -1.19%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.215.58'
$ws.Range("E3").Value = '  -1.74%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
Set-TextValue $ws.Range("D5") '241.54'
$ws.Range("E5").Value = '  -2.20%  '

# Row 6
$ws.Range("E6").Value = '  -0.96%  '

# Row 7
Set-TextValue $ws.Range("D7") '73.21'
$ws.Range("E7").Value = '  -4.00%  '

# Row 8
$ws.Range("E8").Value = '  -0.06%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.605'
$ws.Range("E9").Value = '  -3.35%  '

# Row 10
Set-TextValue $ws.Range("D10") '42.54'
$ws.Range("E10").Value = '  -2.27%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0956'
$ws.Range("E11").Value = '  +0.29%  '

# Row 12
Set-TextValue $ws.Range("D12") '6.99'
$ws.Range("E12").Value = '  -4.73%  '

# Row 13
$ws.Range("E13").Value = '  -0.02%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.552.20'
$ws.Range("E14").Value = '  -1.64%  '

# Row 15
Set-TextValue $ws.Range("D15") '14.18'
$ws.Range("E15").Value = '  -3.35%  '

# Row 16
Set-TextValue $ws.Range("D16") '0.834'
$ws.Range("E16").Value = '  -2.96%  '

# Row 17
Set-TextValue $ws.Range("D17") '2.210.64'
$ws.Range("E17").Value = '  -2.03%  '

# Row 18
Set-TextValue $ws.Range("D18") '41.826.62'
$ws.Range("E18").Value = '  -1.20%  '

# Row 19
$ws.Range("E19").Value = '  +4.47%  '

# Row 20
Set-TextValue $ws.Range("B20") 'Litecoin'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D20") '72.68'
$ws.Range("E20").Value = '  +0.40%  '

# Row 21
Set-TextValue $ws.Range("B21") 'Uniswap'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '6.17'
$ws.Range("E21").Value = '  -0.90%  '

# Row 22
Set-TextValue $ws.Range("D22") '10.86'
$ws.Range("E22").Value = '  +16.62%  '

# Row 23
$ws.Range("E23").Value = '  -0.97%  '

# Row 24
$ws.Range("E24").Value = '  -7.43%  '

# Row 25
Set-TextValue $ws.Range("D25") '11.59'
$ws.Range("E25").Value = '  +0.75%  '

# Row 26
Set-TextValue $ws.Range("D26") '1.00'
$ws.Range("E26").Value = '  +0.33%  '

# Row 27
Set-TextValue $ws.Range("D27") '3.75'
$ws.Range("E27").Value = '  +3.75%  '

# Row 28
$ws.Range("E28").Value = '  -2.00%  '

# Row 29
$ws.Range("E29").Value = '  -1.41%  '

# Row 30
Set-TextValue $ws.Range("D30") '167.68'
$ws.Range("E30").Value = '  -0.76%  '

# Row 31
Set-TextValue $ws.Range("D31") '20.50'
$ws.Range("E31").Value = '  -1.27%  '

# Row 32
Set-TextValue $ws.Range("D32") '5.66'
$ws.Range("E32").Value = '  +4.71%  '

# Row 33
$ws.Range("E33").Value = '  -3.86%  '

# Row 34
Set-TextValue $ws.Range("D34") '30.12'
$ws.Range("E34").Value = '  -3.02%  '

# Row 35
$ws.Range("E35").Value = '  -0.60%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.110'
$ws.Range("E36").Value = '  -9.79%  '

# Row 37
Set-TextValue $ws.Range("D37") '4.24'
$ws.Range("E37").Value = '  -6.12%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.0300'
$ws.Range("E38").Value = '  -5.93%  '

# Row 39
Set-TextValue $ws.Range("D39") '13.78'
$ws.Range("E39").Value = '  -0.24%  '

# Row 40
Set-TextValue $ws.Range("D40") '65.03'
$ws.Range("E40").Value = '  +2.25%  '

# Row 41
$ws.Range("E41").Value = '  -3.50%  '

# Row 42
$ws.Range("E42").Value = '  -3.35%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.197'
$ws.Range("E43").Value = '  -3.72%  '

# Row 44
Set-TextValue $ws.Range("D44") '8.77'
$ws.Range("E44").Value = '  -0.74%  '

# Row 45
Set-TextValue $ws.Range("D45") '104.86'
$ws.Range("E45").Value = '  -4.00%  '

# Row 46
$ws.Range("E46").Value = '  -2.57%  '

# Row 47
$ws.Range("E47").Value = '  +2.12%  '

# Row 48
Set-TextValue $ws.Range("B48") 'ARBITRUM'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D48") '1.12'
$ws.Range("E48").Value = '  -1.35%  '

# Row 49
Set-TextValue $ws.Range("B49") 'TrustWalletToken'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D49") '1.17'
$ws.Range("E49").Value = '  -2.51%  '

# Row 50
$ws.Range("E50").Value = '  +0.03%  '

# Row 51
Set-TextValue $ws.Range("D51") '2.424.29'
$ws.Range("E51").Value = '  -1.79%  '

